# Apply the cryptos-list price/volume update described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel
# (single decimal point, e.g. "217.40") must be pre-formatted as Text so they
# round-trip as the literal string, matching the original inline-string data.
$textCells = @("D5", "D10", "D14", "D15", "D16", "D19", "D21", "D24", "D25", "D26", "D29", "D34", "D37", "D39", "D41", "D42", "D45", "D46", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.239.61"
$ws.Range("E2").Value = "  +1.47%  "
$ws.Range("D3").Value = "1.645.08"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "217.40"
$ws.Range("E5").Value = "  +0.39%  "
$ws.Range("E6").Value = "  +1.84%  "
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("E8").Value = "  +1.37%  "
$ws.Range("E9").Value = "  +1.08%  "
$ws.Range("D10").Value = "20.06"
$ws.Range("E10").Value = "  +1.67%  "
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("D12").Value = "1.875.42"
$ws.Range("E12").Value = "  +0.17%  "
$ws.Range("D13").Value = "1.629.33"
$ws.Range("E13").Value = "  -1.35%  "
$ws.Range("D14").Value = "4.16"
$ws.Range("E14").Value = "  +1.05%  "
$ws.Range("D15").Value = "0.545"
$ws.Range("E15").Value = "  +3.22%  "
$ws.Range("D16").Value = "67.40"
$ws.Range("E16").Value = "  +2.17%  "
$ws.Range("D17").Value = "27.238.84"
$ws.Range("E17").Value = "  +1.36%  "
$ws.Range("E18").Value = "  +1.96%  "
$ws.Range("D19").Value = "219.81"
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("D21").Value = "6.91"
$ws.Range("E21").Value = "  +4.82%  "
$ws.Range("E22").Value = "  +7.41%  "
$ws.Range("E23").Value = "  +0.91%  "
$ws.Range("D24").Value = "9.23"
$ws.Range("E24").Value = "  +0.94%  "
$ws.Range("D25").Value = "147.93"
$ws.Range("E25").Value = "  +1.31%  "
$ws.Range("D26").Value = "7.53"
$ws.Range("E26").Value = "  +1.65%  "
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("D29").Value = "15.80"
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("E30").Value = "  +1.66%  "
$ws.Range("E31").Value = "  +0.64%  "
$ws.Range("E32").Value = "  +1.18%  "
$ws.Range("E33").Value = "  +1.25%  "
$ws.Range("D34").Value = "1.57"
$ws.Range("E34").Value = "  +1.56%  "
$ws.Range("D35").Value = "1.275.87"
$ws.Range("E35").Value = "  +2.52%  "
$ws.Range("E36").Value = "  +0.97%  "
$ws.Range("D37").Value = "0.0177"
$ws.Range("E37").Value = "  +1.84%  "
$ws.Range("E38").Value = "  +4.18%  "
$ws.Range("D39").Value = "0.546"
$ws.Range("E39").Value = "  +1.79%  "
$ws.Range("E40").Value = "  -0.20%  "
$ws.Range("D41").Value = "0.810"
$ws.Range("E41").Value = "  +0.76%  "
$ws.Range("D42").Value = "2.23"
$ws.Range("E42").Value = "  +6.63%  "
$ws.Range("E43").Value = "  -0.46%  "
$ws.Range("D44").Value = "1.785.58"
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").Value = "61.88"
$ws.Range("E45").Value = "  +1.80%  "
$ws.Range("D46").Value = "91.95"
$ws.Range("E46").Value = "  +0.78%  "
$ws.Range("E47").Value = "  +2.38%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "0.0516"
$ws.Range("E48").Value = "  +0.20%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "7.67"
$ws.Range("E49").Value = "  +1.29%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.0976"
$ws.Range("E50").Value = "  +0.61%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "0.405"
$ws.Range("E51").Value = "  +0.03%  "
